$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (English record): SanDisk Cruzer Blade USB drive -> Dell Vostro laptop
$ws.Range("B2").Value = "Vostro"
$ws.Range("C2").Value = "Dell"
$ws.Range("D2").Value = 3568
$ws.Range("E2").Value = "DKS"
$ws.Range("F2").Value = 1.454
$ws.Range("G2").Value = "To take enrollments"

# Row 3 (Arabic record): same change, Arabic translations
$ws.Range("B3").Value = "ستر  "
$ws.Range("C3").Value = "دلّ  "
$ws.Range("D3").Value = 3568
$ws.Range("E3").Value = "DKS"
$ws.Range("F3").Value = 1.454
$ws.Range("G3").Value = "لأخذ التسجيلات"

Write-Output "done"
